$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 271, shifting existing rows 271-323 down to 272-324
$ws.Rows(271).Insert()

# Populate the new row 271 with the latest weekly price record
$ws.Cells.Item(271, 1).Value = 11
$ws.Cells.Item(271, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(271, 3).Value = "Bíobío"
$ws.Cells.Item(271, 4).Value = 45258
$ws.Cells.Item(271, 5).Value = 8
$ws.Cells.Item(271, 6).Value = "Fruta"
$ws.Cells.Item(271, 7).Value = 100108
$ws.Cells.Item(271, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(271, 9).Value = 100108005
$ws.Cells.Item(271, 10).Value = "Piña"
$ws.Cells.Item(271, 11).Value = "Caramelo"
$ws.Cells.Item(271, 12).Value = "Segunda"
$ws.Cells.Item(271, 13).Value = 100
$ws.Cells.Item(271, 14).Value = 21000
$ws.Cells.Item(271, 15).Value = 22000
$ws.Cells.Item(271, 16).Value = 21500
$ws.Cells.Item(271, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(271, 18).Value = "Ecuador"
$ws.Cells.Item(271, 19).Value = 1536
$ws.Cells.Item(271, 20).Value = 14
